$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) to be treated as text so values like
# "1.00", "530.19", "0.0544" etc. are not auto-converted to numbers.
$priceCells = "D2","D3","D4","D5","D6","D8","D9","D10","D12","D14","D15","D16","D17","D19","D21","D23","D24","D25","D26","D27","D28","D29","D30","D33","D34","D37","D38","D40","D42","D43","D44","D45","D47","D48","D49","D51"
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "61.005.35"
$ws.Range("E2").Value = "  +3.90%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.722.54"
$ws.Range("E3").Value = "  +3.41%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.04%  "

# Row 5 - BNB
$ws.Range("D5").Value = "530.19"
$ws.Range("E5").Value = "  +1.84%  "

# Row 6 - Solana
$ws.Range("D6").Value = "148.34"
$ws.Range("E6").Value = "  +2.29%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.18%  "

# Row 8 - XRP
$ws.Range("D8").Value = "0.580"
$ws.Range("E8").Value = "  +0.74%  "

# Row 9 - LidoStakedEther
$ws.Range("D9").Value = "2.744.95"
$ws.Range("E9").Value = "  +3.96%  "

# Row 10 - Toncoin
$ws.Range("D10").Value = "7.11"
$ws.Range("E10").Value = "  +13.36%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  +1.78%  "

# Row 12 - Cardano
$ws.Range("D12").Value = "0.343"
$ws.Range("E12").Value = "  +2.59%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  +3.01%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "3.200.69"
$ws.Range("E14").Value = "  +3.51%  "

# Row 15 - WrappedBTC
$ws.Range("D15").Value = "60.998.33"
$ws.Range("E15").Value = "  +3.92%  "

# Row 16 - Avalanche
$ws.Range("D16").Value = "21.59"
$ws.Range("E16").Value = "  +4.09%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "2.747.66"
$ws.Range("E17").Value = "  +3.92%  "

# Row 18 - ShibaInu
$ws.Range("E18").Value = "  +1.75%  "

# Row 19 - BitcoinCash
$ws.Range("D19").Value = "346.37"
$ws.Range("E19").Value = "  -0.53%  "

# Row 20 - Polkadot
$ws.Range("E20").Value = "  +1.75%  "

# Row 21 - Chainlink
$ws.Range("D21").Value = "10.60"
$ws.Range("E21").Value = "  +3.89%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  +5.16%  "

# Row 23 - Dai
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  +0.03%  "

# Row 24 - Litecoin
$ws.Range("D24").Value = "63.52"
$ws.Range("E24").Value = "  +3.04%  "

# Row 25 - Kaspa
$ws.Range("D25").Value = "0.172"
$ws.Range("E25").Value = "  +4.98%  "

# Row 26 - Polygon
$ws.Range("D26").Value = "0.420"
$ws.Range("E26").Value = "  +1.62%  "

# Row 27 - Binance-PegBSC-USD
$ws.Range("D27").Value = "0.993"
$ws.Range("E27").Value = "  -0.44%  "

# Row 28 - PEPE
$ws.Range("D28").Value = "0.0₃0828"
$ws.Range("E28").Value = "  +3.01%  "

# Row 29 - InternetComputer(DFINITY)
$ws.Range("D29").Value = "7.35"
$ws.Range("E29").Value = "  +4.84%  "

# Row 30 - Aptos
$ws.Range("D30").Value = "6.79"
$ws.Range("E30").Value = "  +8.98%  "

# Row 31 - USDe
$ws.Range("E31").Value = "  -0.12%  "

# Row 32 - PancakeSwap
$ws.Range("E32").Value = "  +2.48%  "

# Row 33 - EthereumClassic
$ws.Range("D33").Value = "19.10"
$ws.Range("E33").Value = "  +1.34%  "

# Row 34 - Monero
$ws.Range("D34").Value = "150.22"
$ws.Range("E34").Value = "  +0.76%  "

# Row 35 - NEARProtocol
$ws.Range("E35").Value = "  +7.17%  "

# Row 36 - ImmutableX
$ws.Range("E36").Value = "  +8.54%  "

# Row 37 - SuiNetwork
$ws.Range("D37").Value = "0.925"
$ws.Range("E37").Value = "  -6.25%  "

# Row 38 - Fetch.AI
$ws.Range("D38").Value = "0.908"
$ws.Range("E38").Value = "  +7.87%  "

# Row 39 - Stacks
$ws.Range("E39").Value = "  +9.08%  "

# Row 40 - OKB
$ws.Range("D40").Value = "37.33"
$ws.Range("E40").Value = "  +2.07%  "

# Row 41 - Filecoin
$ws.Range("E41").Value = "  +2.12%  "

# Row 42 - Mantle
$ws.Range("D42").Value = "0.628"
$ws.Range("E42").Value = "  +4.59%  "

# Row 43 - Bittensor
$ws.Range("D43").Value = "282.29"
$ws.Range("E43").Value = "  +0.86%  "

# Row 44 - EnergySwap
$ws.Range("D44").Value = "20.31"
$ws.Range("E44").Value = "  +3.78%  "

# Row 45 - Stellar
$ws.Range("D45").Value = "0.0989"
$ws.Range("E45").Value = "  +0.55%  "

# Row 46 - FirstDigitalUSD
$ws.Range("E46").Value = "  -0.14%  "

# Row 47 - Maker
$ws.Range("D47").Value = "2.123.49"
$ws.Range("E47").Value = "  +6.80%  "

# Row 48 - was Hedera, now RenderToken (rows 48/49 content swapped)
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "5.00"
$ws.Range("E48").Value = "  +7.76%  "

# Row 49 - was RenderToken, now Hedera
$ws.Range("B49").Value = "Hedera"
$ws.Range("C49").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D49").Value = "0.0544"
$ws.Range("E49").Value = "  +4.39%  "

# Row 50 - WhiteBITCoin
$ws.Range("E50").Value = "  +2.20%  "

# Row 51 - InjectiveProtocol
$ws.Range("D51").Value = "19.50"
$ws.Range("E51").Value = "  +6.09%  "
